$d = $word.ActiveDocument

# The eight "plan" bullet paragraphs (originally indices 2-9, right after
# the "Common plan:" heading paragraph) are replaced, as a block, by six
# bullet paragraphs: two of the old items are dropped, two are reworded,
# and three gain <w:proofErr/> wrapped sub-runs around a word/phrase that
# Word's proofing pass flagged (spell-check on "ForestNode"/"QGraphicsView",
# grammar check on "checking").

$startPara = $d.Paragraphs(2)
$endPara   = $d.Paragraphs(9)
$target = $d.Range($startPara.Range.Start, $endPara.Range.End)

$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$p1 = '<w:p ' + $ns + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Add points</w:t></w:r></w:p>'

$p2 = '<w:p ' + $ns + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Implement </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>ForestNode</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> as public using metaprogramming</w:t></w:r></w:p>'

$p3 = '<w:p ' + $ns + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Clean up repo and make Tree standalone (and maybe other potential libraries)</w:t></w:r></w:p>'

$p4 = '<w:p ' + $ns + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Solve issue with artifacts on </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>QGraphicsView</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'

$p5 = '<w:p ' + $ns + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Add current sources and different types of loads and converters, modify the file schema</w:t></w:r></w:p>'

$p6 = '<w:p ' + $ns + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Implement setting, </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>checking</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> and displaying limitations</w:t></w:r></w:p>'

$body = $p1 + $p2 + $p3 + $p4 + $p5 + $p6

$wordXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document ' + $ns + '><w:body>' + $body + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($wordXml) | Out-Null
